$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 1724901.6
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 6515.0625
$ws.Range("I33").Value = 10343.1
$ws.Range("J33").Value = 135
$ws.Range("K33").Value = 10343.1
$ws.Range("L33").Value = 135
$ws.Range("M33").Value = -10114.1
$ws.Range("N33").Value = -593
# Row 80 (Leve Item ID 12605)
$ws.Range("H80").Value = 247.5
$ws.Range("I80").Value = 191.78572
$ws.Range("K80").Value = 575.35716
$ws.Range("M80").Value = 422.64284
# Row 83 (Leve Item ID 12605)
$ws.Range("H83").Value = 247.5
$ws.Range("I83").Value = 191.78572
$ws.Range("K83").Value = 1726.07148
$ws.Range("M83").Value = 3265.92852
# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 785710.2
$ws.Range("J125").Value = 1471542.6
$ws.Range("L125").Value = 13243883.4
$ws.Range("N125").Value = -13248803.4
# Row 126 (Leve Item ID 34391)
$ws.Range("H126").Value = 13496.875
$ws.Range("J126").Value = 13496.875
$ws.Range("L126").Value = 13496.875
$ws.Range("N126").Value = -23376.875
# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 1135.3846
$ws.Range("J129").Value = 1138.9474
$ws.Range("L129").Value = 3416.8422
$ws.Range("N129").Value = -13416.8422
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1152.1207
$ws.Range("I137").Value = 1010.7317
$ws.Range("J137").Value = 1493.1177
$ws.Range("K137").Value = 3032.1951
$ws.Range("L137").Value = 4479.3531
$ws.Range("M137").Value = -482.1950999999999
$ws.Range("N137").Value = -9579.3531
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2856.0752
$ws.Range("I138").Value = 1365.2642
$ws.Range("J138").Value = 4831.4
$ws.Range("K138").Value = 4095.7926
$ws.Range("L138").Value = 14494.2
$ws.Range("M138").Value = 1044.2074
$ws.Range("N138").Value = -24774.2

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 2514.72
$ws.Range("I32").Value = 2494.8164
$ws.Range("J32").Value = 3490
$ws.Range("K32").Value = 2494.8164
$ws.Range("L32").Value = 3490
$ws.Range("M32").Value = -2207.8164
$ws.Range("N32").Value = -4064
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 4958.7812
$ws.Range("I61").Value = 6148.591
$ws.Range("J61").Value = 2341.2
$ws.Range("K61").Value = 6148.591
$ws.Range("L61").Value = 2341.2
$ws.Range("M61").Value = -5936.591
$ws.Range("N61").Value = -2765.2
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 4958.7812
$ws.Range("I136").Value = 6148.591
$ws.Range("J136").Value = 2341.2
$ws.Range("K136").Value = 18445.773
$ws.Range("L136").Value = 7023.599999999999
$ws.Range("M136").Value = -15895.773
$ws.Range("N136").Value = -12123.6

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1883.5476
$ws.Range("I134").Value = 1339.8462
$ws.Range("J134").Value = 2767.0625
$ws.Range("K134").Value = 4019.5386
$ws.Range("L134").Value = 8301.1875
$ws.Range("M134").Value = -1484.5386
$ws.Range("N134").Value = -13371.1875

$ws = $wb.Worksheets.Item("CRP")
# Row 6 (Leve Item ID 2219)
$ws.Range("H6").Value = 3989146.5
$ws.Range("I6").Value = 6474113
$ws.Range("K6").Value = 6474113
$ws.Range("M6").Value = -6474000
# Row 17 (Leve Item ID 1823)
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 5000
$ws.Range("K17").Value = 5000
$ws.Range("M17").Value = -4826
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2233.6545
$ws.Range("I31").Value = 1441.2
$ws.Range("J31").Value = 3620.45
$ws.Range("K31").Value = 1441.2
$ws.Range("L31").Value = 3620.45
$ws.Range("M31").Value = -1146.2
$ws.Range("N31").Value = -4210.45
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2233.6545
$ws.Range("I34").Value = 1441.2
$ws.Range("J34").Value = 3620.45
$ws.Range("K34").Value = 1441.2
$ws.Range("L34").Value = 3620.45
$ws.Range("M34").Value = -1239.2
$ws.Range("N34").Value = -4024.45
# Row 68 (Leve Item ID 10611)
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0
# Row 71 (Leve Item ID 10611)
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1735.18
$ws.Range("I134").Value = 1767.3948
$ws.Range("J134").Value = 1633.1666
$ws.Range("K134").Value = 5302.1844
$ws.Range("L134").Value = 4899.4998
$ws.Range("M134").Value = -2767.1844
$ws.Range("N134").Value = -9969.4998

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 1033471.9
$ws.Range("I5").Value = 971.625
$ws.Range("J5").Value = 1951249.9
$ws.Range("K5").Value = 2914.875
$ws.Range("L5").Value = 5853749.699999999
$ws.Range("M5").Value = -2802.875
$ws.Range("N5").Value = -5853973.699999999
# Row 17 (Leve Item ID 4640)
$ws.Range("H17").Value = 970
$ws.Range("I17").Value = 281.25
$ws.Range("J17").Value = 1757.1428
$ws.Range("K17").Value = 843.75
$ws.Range("L17").Value = 5271.428400000001
$ws.Range("M17").Value = -674.75
$ws.Range("N17").Value = -5609.428400000001
# Row 36 (Leve Item ID 4732)
$ws.Range("H36").Value = 1800
$ws.Range("I36").Value = 1500
$ws.Range("J36").Value = 1950
$ws.Range("K36").Value = 4500
$ws.Range("L36").Value = 5850
$ws.Range("M36").Value = -4331
$ws.Range("N36").Value = -6188
# Row 39 (Leve Item ID 4712)
$ws.Range("H39").Value = 9730.885
$ws.Range("J39").Value = 6208.3335
$ws.Range("L39").Value = 18625.0005
$ws.Range("N39").Value = -19213.0005
# Row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 13697.5
$ws.Range("J55").Value = 13697.5
$ws.Range("L55").Value = 41092.5
$ws.Range("N55").Value = -41446.5
# Row 75 (Leve Item ID 12863)
$ws.Range("H75").Value = 993.8889
$ws.Range("I75").Value = 900
$ws.Range("J75").Value = 1005.625
$ws.Range("K75").Value = 2700
$ws.Range("L75").Value = 3016.875
$ws.Range("M75").Value = -1702
$ws.Range("N75").Value = -5012.875
# Row 76 (Leve Item ID 12869)
$ws.Range("H76").Value = 2000
$ws.Range("J76").Value = 2000
$ws.Range("L76").Value = 6000
$ws.Range("N76").Value = -6766
# Row 78 (Leve Item ID 12863)
$ws.Range("H78").Value = 993.8889
$ws.Range("I78").Value = 900
$ws.Range("J78").Value = 1005.625
$ws.Range("K78").Value = 8100
$ws.Range("L78").Value = 9050.625
$ws.Range("M78").Value = -3108
$ws.Range("N78").Value = -19034.625
# Row 79 (Leve Item ID 12869)
$ws.Range("H79").Value = 2000
$ws.Range("J79").Value = 2000
$ws.Range("L79").Value = 6000
$ws.Range("N79").Value = -8652
# Row 81 (Leve Item ID 12843)
$ws.Range("H81").Value = 4338
$ws.Range("J81").Value = 4338
$ws.Range("L81").Value = 13014
$ws.Range("N81").Value = -15260
# Row 82 (Leve Item ID 12856)
$ws.Range("H82").Value = 104400.7
$ws.Range("J82").Value = 148284.86
$ws.Range("L82").Value = 444854.58
$ws.Range("N82").Value = -445666.58
# Row 84 (Leve Item ID 12843)
$ws.Range("H84").Value = 4338
$ws.Range("J84").Value = 4338
$ws.Range("L84").Value = 39042
$ws.Range("N84").Value = -50274
# Row 85 (Leve Item ID 12856)
$ws.Range("H85").Value = 104400.7
$ws.Range("J85").Value = 148284.86
$ws.Range("L85").Value = 444854.58
$ws.Range("N85").Value = -447662.58
# Row 87 (Leve Item ID 12864)
$ws.Range("H87").Value = 10237.5
$ws.Range("I87").Value = 5580
$ws.Range("J87").Value = 18000
$ws.Range("K87").Value = 16740
$ws.Range("L87").Value = 54000
$ws.Range("M87").Value = -15492
$ws.Range("N87").Value = -56496
# Row 88 (Leve Item ID 12851)
$ws.Range("H88").Value = 2960
$ws.Range("J88").Value = 2960
$ws.Range("L88").Value = 8880
$ws.Range("N88").Value = -9736
# Row 90 (Leve Item ID 12864)
$ws.Range("H90").Value = 10237.5
$ws.Range("I90").Value = 5580
$ws.Range("J90").Value = 18000
$ws.Range("K90").Value = 50220
$ws.Range("L90").Value = 162000
$ws.Range("M90").Value = -43980
$ws.Range("N90").Value = -174480
# Row 91 (Leve Item ID 12851)
$ws.Range("H91").Value = 2960
$ws.Range("J91").Value = 2960
$ws.Range("L91").Value = 8880
$ws.Range("N91").Value = -11844
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 5222.353
$ws.Range("I131").Value = 827.5
$ws.Range("J131").Value = 5808.3335
$ws.Range("K131").Value = 2482.5
$ws.Range("L131").Value = 17425.0005
$ws.Range("M131").Value = 2557.5
$ws.Range("N131").Value = -27505.0005
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 1033471.9
$ws.Range("I135").Value = 971.625
$ws.Range("J135").Value = 1951249.9
$ws.Range("K135").Value = 8744.625
$ws.Range("L135").Value = 17561249.1
$ws.Range("M135").Value = -6209.625
$ws.Range("N135").Value = -17566319.1

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 1576.3125
$ws.Range("I113").Value = 1585.0834
$ws.Range("K113").Value = 1585.0834
$ws.Range("M113").Value = 584.9166

$ws = $wb.Worksheets.Item("LTW")
# Row 98 (Leve Item ID 18379)
$ws.Range("H98").Value = 40463.75
$ws.Range("J98").Value = 40463.75
$ws.Range("L98").Value = 40463.75
$ws.Range("N98").Value = -46453.75
# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 6190.5806
$ws.Range("I132").Value = 6373.289
$ws.Range("J132").Value = 5706.9414
$ws.Range("K132").Value = 19119.867
$ws.Range("L132").Value = 17120.8242
$ws.Range("M132").Value = -16589.867
$ws.Range("N132").Value = -22180.8242
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 17546508
$ws.Range("I136").Value = 2964.0715
$ws.Range("K136").Value = 8892.2145
$ws.Range("M136").Value = -6342.2145

$ws = $wb.Worksheets.Item("WVR")
# Row 94 (Leve Item ID 18075)
$ws.Range("H94").Value = 29500
$ws.Range("J94").Value = 29500
$ws.Range("L94").Value = 29500
$ws.Range("N94").Value = -31302
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2278.2703
$ws.Range("I132").Value = 2077.1936
$ws.Range("J132").Value = 3317.1667
$ws.Range("K132").Value = 6231.5808
$ws.Range("L132").Value = 9951.500100000001
$ws.Range("M132").Value = -3701.5808
$ws.Range("N132").Value = -15011.5001
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 3581.1738
$ws.Range("I136").Value = 542.70966
$ws.Range("J136").Value = 9860.666999999999
$ws.Range("K136").Value = 1628.12898
$ws.Range("L136").Value = 29582.001
$ws.Range("M136").Value = 921.87102
$ws.Range("N136").Value = -34682.001
